# Updated symbol list on Mon Jan  9 08:19:17 UTC 2023 with GitHub Actions
# Applies updated Price (D), Volume(1h) (E), and Hora (G) values to the
# cryptos worksheet. Values are written as literal text (NumberFormat "@")
# to preserve exact formatting (trailing zeros, percent signs, etc.)
# matching the original inline-string cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.32%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "8"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.61%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "8"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.798"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.80%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "8"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06254"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.46%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "8"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.923"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.76%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "8"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8804"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.59%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "8"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9432"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.41%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1463"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.16%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "8"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05354"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.96%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07341"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "8"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03111"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.84%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "8"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09058"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.09%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "8"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001553"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.85%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "8"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006278"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.69%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "8"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005720"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.75%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "8"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.443"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.10%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "8"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.266"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.99%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "8"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.48%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "8"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3147"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.30%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "8"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1313"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.17%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "8"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.848"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.90%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "8"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04309"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.63%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "8"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001179"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.76%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "8"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004277"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.93%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "8"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.03%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "8"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001692"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.21%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "8"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "8"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "8"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "8"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "8"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "8"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "8"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "8"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "8"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "8"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "8"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "8"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "8"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.58%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "8"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006709"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "62.30%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "8"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.76%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "8"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002145"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.11%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "8"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01300"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.53%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "8"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005108"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.04%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "8"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.04%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "8"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "839.26%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "8"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "8"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.04%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "8"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "8"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "8"
